# Applies the slide-12 (sldId 658) geometry tweaks, removes the stray
# "Slide Number Placeholder" shape, and fixes the "containts" typo on
# slide 16 (sldId 733).

function EmuToPt($emu) {
    # Shape.Left/.Top/.Width/.Height are in points and are stored
    # internally with limited (single-precision-ish) resolution, so a
    # literal EMU/12700 division can truncate down by 1 EMU once
    # re-serialised. A small nudge keeps the round-trip exact.
    return ($emu / 12700.0) + 0.00003
}

$p = $ppt.ActivePresentation

# ---- Slide 12 (sldId 658) -------------------------------------------------
$s = $p.Slides.Item(12)

# Flowchart: Connector 13 (id 14) - taller, shifted up
$sh = $s.Shapes.Item("Flowchart: Connector 13")
$sh.Top    = EmuToPt 358525
$sh.Height = EmuToPt 2382842

# Flowchart: Connector 29 (id 30) - taller, shifted up
$sh = $s.Shapes.Item("Flowchart: Connector 29")
$sh.Top    = EmuToPt 766365
$sh.Height = EmuToPt 2178002

# TextBox 30 (id 31) - shifted up
$sh = $s.Shapes.Item("TextBox 30")
$sh.Top = EmuToPt 98013

# Flowchart: Connector 31 (id 32) - taller, shifted up
$sh = $s.Shapes.Item("Flowchart: Connector 31")
$sh.Top    = EmuToPt 2762904
$sh.Height = EmuToPt 3725960

# Flowchart: Connector 32 (id 33) - taller, shifted up
$sh = $s.Shapes.Item("Flowchart: Connector 32")
$sh.Top    = EmuToPt 3565599
$sh.Height = EmuToPt 2178002

# Remove the slide-number placeholder entirely. Deleting a placeholder
# once re-materialises an empty stub (mirrors the layout), so delete twice.
$sh = $s.Shapes.Item("Slide Number Placeholder 3")
$sh.Delete()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -like "Slide Number Placeholder*") {
        $cand.Delete()
        break
    }
}

# ---- Slide 16 (sldId 733) --------------------------------------------------
# Fix the "containts" typo: the paragraph was split across three runs
# ("It " / "containts" / " the reference to the object "); replace the
# whole paragraph range (Characters(start,length) keeps run-merging
# behaviour clean, unlike Paragraphs().Text which only patches inside
# existing run boundaries) with a single corrected run.
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item("Content Placeholder 2")
$tr = $sh16.TextFrame.TextRange
$para4 = $tr.Paragraphs(4, 1)
$fullPara = $tr.Characters($para4.Start, $para4.Length)
$fullPara.Text = "It contains the reference to the object "
